# Weekly refresh: insert the newest "Fruta / hortaliza, semanal" price group
# (3 rows: Especial / Primera / Segunda) at the top of the Palta data block,
# pushing all existing rows down by 3 (dimension grows from T665 to T668).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 blank rows above row 635 (existing data - previously rows 635:665 -
# shifts down to 638:668).
$ws.Rows.Item(635).Resize(3).Insert()

# Columns: A Mercado ID, B Mercado, C Región, D Fecha, E Codreg, F Tipo,
#          G Producto ID, H Producto, I Categoría ID, J Categoría, K Variedad,
#          L Calidad, M Volumen, N Precio mínimo, O Precio máximo,
#          P Precio promedio ponderado, Q Unidad de comercialización,
#          R Origen, S Precio $/Kg, T Kg / unidad

$rows = @(
    @{ Row=635; L="Especial"; M=300; N=2600; O=2700; P=2650; S=2650 },
    @{ Row=636; L="Primera";  M=240; N=2300; O=2400; P=2350; S=2350 },
    @{ Row=637; L="Segunda";  M=200; N=2000; O=2100; P=2050; S=2050 }
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value  = 8
    $ws.Cells.Item($row, 2).Value  = "Terminal La Palmera de La Serena"
    $ws.Cells.Item($row, 3).Value  = "Coquimbo"
    $ws.Cells.Item($row, 4).Value  = 44610
    $ws.Cells.Item($row, 5).Value  = 4
    $ws.Cells.Item($row, 6).Value  = "Fruta"
    $ws.Cells.Item($row, 7).Value  = 100106
    $ws.Cells.Item($row, 8).Value  = "Oleaginosos"
    $ws.Cells.Item($row, 9).Value  = 100106002
    $ws.Cells.Item($row, 10).Value = "Palta"
    $ws.Cells.Item($row, 11).Value = "Hass"
    $ws.Cells.Item($row, 12).Value = $r.L
    $ws.Cells.Item($row, 13).Value = $r.M
    $ws.Cells.Item($row, 14).Value = $r.N
    $ws.Cells.Item($row, 15).Value = $r.O
    $ws.Cells.Item($row, 16).Value = $r.P
    $ws.Cells.Item($row, 17).Value = "`$/kilo (en caja de 17 kilos)"
    $ws.Cells.Item($row, 18).Value = "Provincia de Limarí"
    $ws.Cells.Item($row, 19).Value = $r.S
    $ws.Cells.Item($row, 20).Value = 1
}
